# Progress log update: append a "Week 9" block (rows 40-47) after the
# existing "Week 8" block, following the same layout/formatting pattern
# used by every earlier week in the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 40: "Week 9" section header (bold label, like A35 "Week 8") ---
$ws.Range("A40").Value = "Week 9"
$ws.Range("A40").Font.Bold = $true

# --- Rows 41-46: daily entries ---
$days   = @(42263, 42264, 42265, 42266, 42267, 42268)
$froms  = @(0.45833333333333331, 0.875, 0.5, 0.54166666666666663, 0.95833333333333337, 0.45833333333333331)
$tos    = @(0.20833333333333334, 0.083333333333333329, 0.041666666666666664, 0.58333333333333337, 0.083333333333333329, 0)
$breaks = @(7, 1, 4, 0, 0, 6)

for ($i = 0; $i -lt 6; $i++) {
    $r = 41 + $i
    $ws.Cells.Item($r, 1).Value = $days[$i]
    $ws.Cells.Item($r, 2).Value = $froms[$i]
    $ws.Cells.Item($r, 3).Value = $tos[$i]
    $ws.Cells.Item($r, 4).Value = $breaks[$i]
}

# Date / time formats matching the rest of the log (columns A/B/C)
$ws.Range("A41:A46").NumberFormat = "[`$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("B41:C46").NumberFormat = "h:mm AM/PM"

# Duration column: one shared formula across the whole week, same pattern
# used by every preceding week block.
$ws.Range("E41:E46").Formula = "=MOD(C41-B41,1)*24-D41"
$ws.Range("E41:E46").NumberFormat = "0.00"

# --- Row 47: week total ---
$ws.Range("D47").Value = "Total"
$ws.Range("D47").Font.Bold = $true
$ws.Range("E47").Formula = "=SUM(E41:E46)"
$ws.Range("E47").Font.Bold = $true
$ws.Range("E47").NumberFormat = "0.00"

# Reflect where the user was working when they saved.
$excel.ActiveWindow.ScrollRow = 25
$ws.Range("E48").Select()

$wb.Save()
